# Update the metrics worksheet so that every model row (2 through 26)
# reflects the new trained model's metrics (all rows now share identical
# values, per the commit "atualizado todo o treinamento para o novo lm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (columns B..Q), identical for every data row (2-26).
# Values are parsed via [double] casts because scientific-notation
# literals (e.g. 3.04e-05) aren't accepted directly by the script parser.
$newValues = @(
    [double]"0.9999674344715328",    # B r2
    [double]"0.9989400190852257",    # C r2_sup
    [double]"0.9999717472652621",    # D r2_test
    [double]"0.9999189651314859",    # E r2_val
    [double]"0.9999480524142762",    # F r2_vt
    [double]"3.039849136261971e-05", # G mse
    [double]"0.0009894456561571146", # H mse_sup
    [double]"3.759280564295439e-05", # I mse_test
    [double]"8.800505942399128e-05", # J mse_val
    [double]"6.279893253347284e-05", # K mse_vt
    [double]"0.000349036811621329",  # L mape
    [double]"0.005513482689065026",  # M rmse
    [double]"1.000060120975632",     # N r2_adj
    [double]"0.005748203061621094",  # O rsd
    [double]"94.80223515402628",     # P aic
    [double]"139.9006406741497"      # Q bic
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
